$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-29 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-30 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("60×80=4800", $true, $false, $false, $false, $false, $true, 1, $false, "27×69=1863", 2) | Out-Null
$d.Content.Find.Execute("76×85=6460", $true, $false, $false, $false, $false, $true, 1, $false, "20×97=1940", 2) | Out-Null
$d.Content.Find.Execute("42×29=1218", $true, $false, $false, $false, $false, $true, 1, $false, "78×87=6786", 2) | Out-Null
$d.Content.Find.Execute("26×55=1430", $true, $false, $false, $false, $false, $true, 1, $false, "62×67=4154", 2) | Out-Null
$d.Content.Find.Execute("44×90=3960", $true, $false, $false, $false, $false, $true, 1, $false, "88×19=1672", 2) | Out-Null
$d.Content.Find.Execute("48×39=1872", $true, $false, $false, $false, $false, $true, 1, $false, "80×17=1360", 2) | Out-Null
$d.Content.Find.Execute("45×52=2340", $true, $false, $false, $false, $false, $true, 1, $false, "79×74=5846", 2) | Out-Null
$d.Content.Find.Execute("40×78=3120", $true, $false, $false, $false, $false, $true, 1, $false, "21×29=609", 2) | Out-Null
$d.Content.Find.Execute("87×97=8439", $true, $false, $false, $false, $false, $true, 1, $false, "21×85=1785", 2) | Out-Null
$d.Content.Find.Execute("44×43=1892", $true, $false, $false, $false, $false, $true, 1, $false, "19×33=627", 2) | Out-Null
$d.Content.Find.Execute("87×44=3828", $true, $false, $false, $false, $false, $true, 1, $false, "35×14=490", 2) | Out-Null
$d.Content.Find.Execute("87×99=8613", $true, $false, $false, $false, $false, $true, 1, $false, "51×97=4947", 2) | Out-Null
$d.Content.Find.Execute("35×67=2345", $true, $false, $false, $false, $false, $true, 1, $false, "40×92=3680", 2) | Out-Null
$d.Content.Find.Execute("69×19=1311", $true, $false, $false, $false, $false, $true, 1, $false, "20×29=580", 2) | Out-Null
$d.Content.Find.Execute("39×43=1677", $true, $false, $false, $false, $false, $true, 1, $false, "55×34=1870", 2) | Out-Null
$d.Content.Find.Execute("43×86=3698", $true, $false, $false, $false, $false, $true, 1, $false, "81×12=972", 2) | Out-Null
$d.Content.Find.Execute("99×30=2970", $true, $false, $false, $false, $false, $true, 1, $false, "50×72=3600", 2) | Out-Null
$d.Content.Find.Execute("96×45=4320", $true, $false, $false, $false, $false, $true, 1, $false, "51×25=1275", 2) | Out-Null
$d.Content.Find.Execute("65×78=5070", $true, $false, $false, $false, $false, $true, 1, $false, "75×52=3900", 2) | Out-Null
$d.Content.Find.Execute("83×87=7221", $true, $false, $false, $false, $false, $true, 1, $false, "88×26=2288", 2) | Out-Null
$d.Content.Find.Execute("38×36=1368", $true, $false, $false, $false, $false, $true, 1, $false, "53×17=901", 2) | Out-Null
$d.Content.Find.Execute("19×64=1216", $true, $false, $false, $false, $false, $true, 1, $false, "27×51=1377", 2) | Out-Null
$d.Content.Find.Execute("43×43=1849", $true, $false, $false, $false, $false, $true, 1, $false, "27×99=2673", 2) | Out-Null
$d.Content.Find.Execute("68×61=4148", $true, $false, $false, $false, $false, $true, 1, $false, "60×11=660", 2) | Out-Null
$d.Content.Find.Execute("69×78=5382", $true, $false, $false, $false, $false, $true, 1, $false, "79×40=3160", 2) | Out-Null
